$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 37411.75
$ws.Cells.Item(3, 10).Value = 37411.75
$ws.Cells.Item(3, 12).Value = 37411.75
$ws.Cells.Item(3, 14).Value = -37639.75
$ws.Cells.Item(93, 8).Value = 44862.332
$ws.Cells.Item(93, 10).Value = 44862.332
$ws.Cells.Item(93, 12).Value = 44862.332
$ws.Cells.Item(93, 14).Value = -49854.332
$ws.Cells.Item(95, 8).Value = 37996.668
$ws.Cells.Item(95, 10).Value = 37996.668
$ws.Cells.Item(95, 12).Value = 37996.668
$ws.Cells.Item(95, 14).Value = -43488.668
$ws.Cells.Item(102, 8).Value = 37411.75
$ws.Cells.Item(102, 10).Value = 37411.75
$ws.Cells.Item(102, 12).Value = 37411.75
$ws.Cells.Item(102, 14).Value = -43901.75
$ws.Cells.Item(105, 8).Value = 38664
$ws.Cells.Item(105, 10).Value = 38664
$ws.Cells.Item(105, 12).Value = 38664
$ws.Cells.Item(105, 14).Value = -45652

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 34696
$ws.Cells.Item(24, 10).Value = 34696
$ws.Cells.Item(24, 12).Value = 34696
$ws.Cells.Item(24, 14).Value = -35444
$ws.Cells.Item(32, 8).Value = 10131.138
$ws.Cells.Item(32, 9).Value = 7607.5454
$ws.Cells.Item(32, 10).Value = 29562.8
$ws.Cells.Item(32, 11).Value = 7607.5454
$ws.Cells.Item(32, 12).Value = 29562.8
$ws.Cells.Item(32, 13).Value = -7320.5454
$ws.Cells.Item(32, 14).Value = -30136.8
$ws.Cells.Item(95, 8).Value = 38066
$ws.Cells.Item(95, 10).Value = 38066
$ws.Cells.Item(95, 12).Value = 38066
$ws.Cells.Item(95, 14).Value = -43558
$ws.Cells.Item(100, 8).Value = 34696
$ws.Cells.Item(100, 10).Value = 34696
$ws.Cells.Item(100, 12).Value = 34696
$ws.Cells.Item(100, 14).Value = -36860
$ws.Cells.Item(101, 8).Value = 48528
$ws.Cells.Item(101, 10).Value = 48528
$ws.Cells.Item(101, 12).Value = 48528
$ws.Cells.Item(101, 14).Value = -55018
$ws.Cells.Item(103, 8).Value = 39354
$ws.Cells.Item(103, 10).Value = 39354
$ws.Cells.Item(103, 12).Value = 39354
$ws.Cells.Item(103, 14).Value = -41698
$ws.Cells.Item(104, 8).Value = 27408.75
$ws.Cells.Item(104, 10).Value = 27408.75
$ws.Cells.Item(104, 12).Value = 27408.75
$ws.Cells.Item(104, 14).Value = -34396.75
$ws.Cells.Item(105, 8).Value = 43453
$ws.Cells.Item(105, 10).Value = 43453
$ws.Cells.Item(105, 12).Value = 43453
$ws.Cells.Item(105, 14).Value = -50441
$ws.Cells.Item(106, 8).Value = 46244
$ws.Cells.Item(106, 10).Value = 46244
$ws.Cells.Item(106, 12).Value = 46244
$ws.Cells.Item(106, 14).Value = -48768
$ws.Cells.Item(132, 8).Value = 11629508
$ws.Cells.Item(132, 9).Value = 20834002
$ws.Cells.Item(132, 10).Value = 2778.2632
$ws.Cells.Item(132, 11).Value = 62502006
$ws.Cells.Item(132, 12).Value = 8334.7896
$ws.Cells.Item(132, 13).Value = -62499476
$ws.Cells.Item(132, 14).Value = -13394.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value = 28397.2
$ws.Cells.Item(92, 10).Value = 28397.2
$ws.Cells.Item(92, 12).Value = 28397.2
$ws.Cells.Item(92, 14).Value = -33389.2
$ws.Cells.Item(95, 8).Value = 43996
$ws.Cells.Item(95, 10).Value = 43996
$ws.Cells.Item(95, 12).Value = 43996
$ws.Cells.Item(95, 14).Value = -49488
$ws.Cells.Item(96, 8).Value = 13727.333
$ws.Cells.Item(96, 9).Value = 1682
$ws.Cells.Item(96, 10).Value = 19750
$ws.Cells.Item(96, 11).Value = 1682
$ws.Cells.Item(96, 12).Value = 19750
$ws.Cells.Item(96, 13).Value = 1064
$ws.Cells.Item(96, 14).Value = -25242
$ws.Cells.Item(100, 8).Value = 45192
$ws.Cells.Item(100, 10).Value = 45192
$ws.Cells.Item(100, 12).Value = 45192
$ws.Cells.Item(100, 14).Value = -47356
$ws.Cells.Item(103, 8).Value = 74259.60000000001
$ws.Cells.Item(103, 10).Value = 74259.60000000001
$ws.Cells.Item(103, 12).Value = 74259.60000000001
$ws.Cells.Item(103, 14).Value = -76603.60000000001
$ws.Cells.Item(106, 8).Value = 32556.5
$ws.Cells.Item(106, 10).Value = 32556.5
$ws.Cells.Item(106, 12).Value = 32556.5
$ws.Cells.Item(106, 14).Value = -35080.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 39047.668
$ws.Cells.Item(28, 10).Value = 39047.668
$ws.Cells.Item(28, 12).Value = 39047.668
$ws.Cells.Item(28, 14).Value = -39537.668
$ws.Cells.Item(31, 8).Value = 148857.38
$ws.Cells.Item(31, 9).Value = 1804.6818
$ws.Cells.Item(31, 10).Value = 191992.84
$ws.Cells.Item(31, 11).Value = 1804.6818
$ws.Cells.Item(31, 12).Value = 191992.84
$ws.Cells.Item(31, 13).Value = -1509.6818
$ws.Cells.Item(31, 14).Value = -192582.84
$ws.Cells.Item(34, 8).Value = 148857.38
$ws.Cells.Item(34, 9).Value = 1804.6818
$ws.Cells.Item(34, 10).Value = 191992.84
$ws.Cells.Item(34, 11).Value = 1804.6818
$ws.Cells.Item(34, 12).Value = 191992.84
$ws.Cells.Item(34, 13).Value = -1602.6818
$ws.Cells.Item(34, 14).Value = -192396.84
$ws.Cells.Item(43, 8).Value = 21382.834
$ws.Cells.Item(43, 10).Value = 21382.834
$ws.Cells.Item(43, 12).Value = 21382.834
$ws.Cells.Item(43, 14).Value = -21750.834
$ws.Cells.Item(92, 8).Value = 35776.555
$ws.Cells.Item(92, 10).Value = 35776.555
$ws.Cells.Item(92, 12).Value = 35776.555
$ws.Cells.Item(92, 14).Value = -40768.555
$ws.Cells.Item(96, 8).Value = 36786.3
$ws.Cells.Item(96, 10).Value = 36786.3
$ws.Cells.Item(96, 12).Value = 36786.3
$ws.Cells.Item(96, 14).Value = -42278.3
$ws.Cells.Item(101, 8).Value = 21382.834
$ws.Cells.Item(101, 10).Value = 21382.834
$ws.Cells.Item(101, 12).Value = 21382.834
$ws.Cells.Item(101, 14).Value = -27872.834
$ws.Cells.Item(106, 8).Value = 40633.5
$ws.Cells.Item(106, 10).Value = 40633.5
$ws.Cells.Item(106, 12).Value = 40633.5
$ws.Cells.Item(106, 14).Value = -43157.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(101, 8).Value = 49553
$ws.Cells.Item(101, 10).Value = 49553
$ws.Cells.Item(101, 12).Value = 49553
$ws.Cells.Item(101, 14).Value = -56043
$ws.Cells.Item(104, 8).Value = 43570.2
$ws.Cells.Item(104, 10).Value = 43570.2
$ws.Cells.Item(104, 12).Value = 43570.2
$ws.Cells.Item(104, 14).Value = -50558.2
$ws.Cells.Item(118, 8).Value = 38306
$ws.Cells.Item(118, 10).Value = 38306
$ws.Cells.Item(118, 12).Value = 38306
$ws.Cells.Item(118, 14).Value = -41620

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(95, 8).Value = 30368.8
$ws.Cells.Item(95, 10).Value = 30368.8
$ws.Cells.Item(95, 12).Value = 30368.8
$ws.Cells.Item(95, 14).Value = -35860.8
$ws.Cells.Item(97, 8).Value = 29397
$ws.Cells.Item(97, 10).Value = 29397
$ws.Cells.Item(97, 12).Value = 29397
$ws.Cells.Item(97, 14).Value = -31379
$ws.Cells.Item(105, 8).Value = 34722.715
$ws.Cells.Item(105, 10).Value = 34722.715
$ws.Cells.Item(105, 12).Value = 34722.715
$ws.Cells.Item(105, 14).Value = -41710.715
$ws.Cells.Item(106, 8).Value = 36090.5
$ws.Cells.Item(106, 10).Value = 36090.5
$ws.Cells.Item(106, 12).Value = 36090.5
$ws.Cells.Item(106, 14).Value = -38614.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(94, 8).Value = 22765
$ws.Cells.Item(94, 10).Value = 22765
$ws.Cells.Item(94, 12).Value = 22765
$ws.Cells.Item(94, 14).Value = -24567
$ws.Cells.Item(95, 8).Value = 34580
$ws.Cells.Item(95, 10).Value = 34580
$ws.Cells.Item(95, 12).Value = 34580
$ws.Cells.Item(95, 14).Value = -40072
$ws.Cells.Item(97, 8).Value = 32152
$ws.Cells.Item(97, 10).Value = 32152
$ws.Cells.Item(97, 12).Value = 32152
$ws.Cells.Item(97, 14).Value = -34134
$ws.Cells.Item(98, 8).Value = 45000
$ws.Cells.Item(98, 10).Value = 45000
$ws.Cells.Item(98, 12).Value = 45000
$ws.Cells.Item(98, 14).Value = -50990
$ws.Cells.Item(99, 8).Value = 36532.363
$ws.Cells.Item(99, 10).Value = 35554.668
$ws.Cells.Item(99, 12).Value = 35554.668
$ws.Cells.Item(99, 14).Value = -41544.668
$ws.Cells.Item(103, 8).Value = 35770
$ws.Cells.Item(103, 10).Value = 35770
$ws.Cells.Item(103, 12).Value = 35770
$ws.Cells.Item(103, 14).Value = -38114
$ws.Cells.Item(104, 8).Value = 46353
$ws.Cells.Item(104, 10).Value = 46353
$ws.Cells.Item(104, 12).Value = 46353
$ws.Cells.Item(104, 14).Value = -53341
$ws.Cells.Item(105, 8).Value = 39740.668
$ws.Cells.Item(105, 10).Value = 39740.668
$ws.Cells.Item(105, 12).Value = 39740.668
$ws.Cells.Item(105, 14).Value = -46728.668
